$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("Preprint", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EarthArXiv, Preprint", 2)

Write-Output "Found: $found"
